$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42:113 down to 43:114
$ws.Rows("42").Insert()

# Populate the newly inserted row 42 with the new record's data
$ws.Range("A42").Value2 = 5
$ws.Range("B42").Value2 = "Macroferia Regional de Talca"
$ws.Range("C42").Value2 = "Maule"
$ws.Range("D42").Value2 = 44880
$ws.Range("E42").Value2 = 7
$ws.Range("F42").Value2 = 100112022
$ws.Range("G42").Value2 = "Arveja Verde"
$ws.Range("H42").Value2 = "Sin especificar"
$ws.Range("I42").Value2 = "Primera"
$ws.Range("J42").Value2 = 500
$ws.Range("K42").Value2 = 17000
$ws.Range("L42").Value2 = 17000
$ws.Range("M42").Value2 = 17000
$ws.Range("N42").Value2 = "$/saco 25 kilos"
$ws.Range("O42").Value2 = "Región del Maule"
$ws.Range("P42").Value2 = 680
$ws.Range("Q42").Value2 = 25
$ws.Range("R42").Value2 = "Hortaliza"
